# Applies the "Reword some stuff in the report" edit described by the
# supplied diff. The large block of the diff that appears to "move" the
# Target Audience / Why-a-computer-is-suitable / Research / Features and
# limitations sections around is an artefact of a single `_GoBack`
# bookmark relocating (which forces Word to renumber every other
# bookmark id that follows it) -- the actual paragraph text in that
# region is unchanged. The genuine content edits are only:
#
#   1. "...trade-off of easily being able to read..."
#        -> "...trade-off of a human (me in particular) being able to
#            easily read..."
#   2. "GPU: [TODO] (Research into whether WPF can be run without a GPU.
#        I know it uses DirectX 9)"
#        -> "GPU: Integrated graphics card, or better."
#   3. "...Certain features of the game (explained below), and small
#        parts..." -> "...Certain features of the game, and small
#        parts..." (the "(explained below)" aside is dropped)
#   4. The `_GoBack` bookmark (Word's "last edit position" marker) moves
#      from the first paragraph to right after "Certain features of the
#      game" -- i.e. right where edit #3 just happened, which is exactly
#      what real Word does after you type/delete there.
#
# Plus the Table of Contents page-number caches for the two TOC rows
# that fall after the newly-lengthened "Analysis" section shift from
# "3" to "4" (pure field-result cache, recomputed by Word on layout).

$d = $word.ActiveDocument

# --- 1. "trade-off of ... being able to easily read" rewording ---------
$d.Content.Find.Execute(
    "there is a trade-off of easily being able to read and debug the data",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "there is a trade-off of a human (me in particular) being able to easily read and debug the data",
    2) | Out-Null

# --- 2. GPU requirement rewording ---------------------------------------
$d.Content.Find.Execute(
    "GPU: [TODO] (Research into whether WPF can be run without a GPU. I know it uses DirectX 9)",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "GPU: Integrated graphics card, or better.",
    2) | Out-Null

# --- 3. Drop the "(explained below)" aside ------------------------------
$d.Content.Find.Execute(
    "Certain features of the game (explained below), and small parts of the code",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Certain features of the game, and small parts of the code",
    2) | Out-Null

# --- 4. Move the _GoBack bookmark to the new last-edit location --------
# Adding a bookmark named "_GoBack" replaces/moves any existing one of
# that name, mirroring Word's single "last edit" bookmark behaviour.
$goBackRange = $d.Content
$goBackRange.Find.Execute("Certain features of the game") | Out-Null
$goBackRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

# --- 5. Refresh the two affected TOC page-number caches -----------------
$d.Content.Find.Execute(
    "Test Data for beta testing",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$toc1 = $d.Content
$toc1.Find.Execute("Test Data for beta testing") | Out-Null
$toc1.Collapse(0)

$d.Content.Find.Execute(
    "Development" ) | Out-Null
